$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Solar capacity values for 2022 (row 24) and 2024 (row 26)
$ws.Range("E24").Value = 22.36
$ws.Range("E26").Value = 55.42
